$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2570.5715
$ws.Range("I76").Value = 1759
$ws.Range("K76").Value = 1759
$ws.Range("M76").Value = -1444

$ws.Range("H79").Value = 2570.5715
$ws.Range("I79").Value = 1759
$ws.Range("K79").Value = 1759
$ws.Range("M79").Value = -667

$ws.Range("H92").Value = 1736716.1
$ws.Range("I92").Value = 781927.6
$ws.Range("K92").Value = 781927.6
$ws.Range("M92").Value = -780679.6

$ws.Range("H127").Value = 2000
$ws.Range("I127").Value = 2000
$ws.Range("K127").Value = 6000
$ws.Range("M127").Value = -1040

$ws.Range("H135").Value = 607.4
$ws.Range("I135").Value = 386.54544
$ws.Range("J135").Value = 1214.75
$ws.Range("K135").Value = 3478.90896
$ws.Range("L135").Value = 10932.75
$ws.Range("M135").Value = -943.9089599999998
$ws.Range("N135").Value = -16002.75

$ws.Range("H137").Value = 1010.86365
$ws.Range("I137").Value = 950
$ws.Range("K137").Value = 2850
$ws.Range("M137").Value = -300

$ws.Range("H138").Value = 4036.0579
$ws.Range("I138").Value = 5845.0435
$ws.Range("J138").Value = 3131.5652
$ws.Range("K138").Value = 17535.1305
$ws.Range("L138").Value = 9394.695599999999
$ws.Range("M138").Value = -12395.1305
$ws.Range("N138").Value = -19674.6956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 413.33334
$ws.Range("J5").Value = 102
$ws.Range("L5").Value = 102
$ws.Range("N5").Value = -326

$ws.Range("H45").Value = 104026.8
$ws.Range("I45").Value = 129033.5
$ws.Range("K45").Value = 129033.5
$ws.Range("M45").Value = -128656.5

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

$ws.Range("H110").Value = 1501.6875
$ws.Range("J110").Value = 1664.3334
$ws.Range("L110").Value = 1664.3334
$ws.Range("N110").Value = -5754.3334

$ws.Range("H122").Value = 1055
$ws.Range("I122").Value = 110
$ws.Range("K122").Value = 330
$ws.Range("M122").Value = 2120

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 413.33334
$ws.Range("J4").Value = 102
$ws.Range("L4").Value = 102
$ws.Range("N4").Value = -332

$ws.Range("H92").Value = 19901
$ws.Range("J92").Value = 19901
$ws.Range("L92").Value = 19901
$ws.Range("N92").Value = -24893

$ws.Range("H99").Value = 6804
$ws.Range("I99").Value = 7714.8
$ws.Range("K99").Value = 7714.8
$ws.Range("M99").Value = -6216.8

$ws.Range("H105").Value = 4993.3794
$ws.Range("I105").Value = 5837.95
$ws.Range("K105").Value = 5837.95
$ws.Range("M105").Value = -4090.95

$ws.Range("H134").Value = 1738.4773
$ws.Range("I134").Value = 1535.8649
$ws.Range("K134").Value = 4607.5947
$ws.Range("M134").Value = -2072.5947

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2734.5
$ws.Range("J11").Value = 2734.5
$ws.Range("L11").Value = 2734.5
$ws.Range("N11").Value = -3014.5

$ws.Range("H16").Value = 10205563
$ws.Range("I16").Value = 11905578
$ws.Range("J16").Value = 5473.5
$ws.Range("K16").Value = 11905578
$ws.Range("L16").Value = 5473.5
$ws.Range("M16").Value = -11905291
$ws.Range("N16").Value = -6047.5

$ws.Range("H17").Value = 9996.333000000001
$ws.Range("J17").Value = 9994.5
$ws.Range("L17").Value = 9994.5
$ws.Range("N17").Value = -10342.5

$ws.Range("H31").Value = 2624.375
$ws.Range("I31").Value = 2847.375
$ws.Range("J31").Value = 2475.7083
$ws.Range("K31").Value = 2847.375
$ws.Range("L31").Value = 2475.7083
$ws.Range("M31").Value = -2552.375
$ws.Range("N31").Value = -3065.7083

$ws.Range("H34").Value = 2624.375
$ws.Range("I34").Value = 2847.375
$ws.Range("J34").Value = 2475.7083
$ws.Range("K34").Value = 2847.375
$ws.Range("L34").Value = 2475.7083
$ws.Range("M34").Value = -2645.375
$ws.Range("N34").Value = -2879.7083

$ws.Range("H113").Value = 10205563
$ws.Range("I113").Value = 11905578
$ws.Range("J113").Value = 5473.5
$ws.Range("K113").Value = 11905578
$ws.Range("L113").Value = 5473.5
$ws.Range("M113").Value = -11903408
$ws.Range("N113").Value = -9813.5

$ws.Range("H132").Value = 47621708
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 47621708
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 142865124
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -142870184

$ws.Range("H134").Value = 1560.5238
$ws.Range("I134").Value = 1220.1111
$ws.Range("J134").Value = 3603
$ws.Range("K134").Value = 3660.3333
$ws.Range("L134").Value = 10809
$ws.Range("M134").Value = -1125.3333
$ws.Range("N134").Value = -15879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 15391
$ws.Range("J106").Value = 17488.75
$ws.Range("L106").Value = 52466.25
$ws.Range("N106").Value = -54358.25

$ws.Range("H140").Value = 2072.9167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7020319
$ws.Range("I80").Value = 88848.50999999999
$ws.Range("K80").Value = 88848.50999999999
$ws.Range("M80").Value = -87850.50999999999

$ws.Range("H83").Value = 7020319
$ws.Range("I83").Value = 88848.50999999999
$ws.Range("K83").Value = 444242.55
$ws.Range("M83").Value = -439250.55

$ws.Range("H132").Value = 9617013
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H140").Value = 114652.25
$ws.Range("J140").Value = 129633.336
$ws.Range("L140").Value = 129633.336
$ws.Range("N140").Value = -139993.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1659.3235
$ws.Range("J61").Value = 2030.9
$ws.Range("L61").Value = 2030.9
$ws.Range("N61").Value = -2434.9

$ws.Range("H113").Value = 1659.3235
$ws.Range("J113").Value = 2030.9
$ws.Range("L113").Value = 2030.9
$ws.Range("N113").Value = -6370.9

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 18149.05
$ws.Range("I136").Value = 51874.6
$ws.Range("K136").Value = 155623.8
$ws.Range("M136").Value = -153073.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 64000
$ws.Range("J68").Value = 64000
$ws.Range("L68").Value = 64000
$ws.Range("N68").Value = -65622

$ws.Range("H71").Value = 64000
$ws.Range("J71").Value = 64000
$ws.Range("L71").Value = 192000
$ws.Range("N71").Value = -200112

$ws.Range("H107").Value = 319
$ws.Range("I107").Value = 309.7143
$ws.Range("K107").Value = 929.1428999999999
$ws.Range("M107").Value = 990.8571000000001

$ws.Range("H122").Value = 2664.125
$ws.Range("I122").Value = 2048.1538
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 6144.4614
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -3694.4614
$ws.Range("N122").Value = -20900.0005

$ws.Range("H136").Value = 43039.168
$ws.Range("I136").Value = 56503.277
$ws.Range("J136").Value = 2646.8333
$ws.Range("K136").Value = 169509.831
$ws.Range("L136").Value = 7940.499899999999
$ws.Range("M136").Value = -166959.831
$ws.Range("N136").Value = -13040.4999
